$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for affected rows.
# D-column values are written with a leading apostrophe so Excel
# stores them as literal text (preserving trailing zeros / dotted
# thousands formatting) instead of coercing them to numbers; the
# style is then reset to "Normal" so no stray quote-prefix format
# is left on the cell (matches the original unformatted cells).
$ws.Range("D2").Value = "'56.544.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.80%  "
$ws.Range("D3").Value = "'2.372.56"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.85%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'505.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -5.12%  "
$ws.Range("E6").Value = "  -2.59%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("E8").Value = "  -2.72%  "
$ws.Range("D9").Value = "'2.392.70"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.37%  "
$ws.Range("D10").Value = "'0.0964"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").Value = "'0.151"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.36%  "
$ws.Range("D12").Value = "'0.323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.89%  "
$ws.Range("D13").Value = "'4.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -10.43%  "
$ws.Range("D14").Value = "'2.795.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.74%  "
$ws.Range("D15").Value = "'56.388.56"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.80%  "
$ws.Range("D16").Value = "'21.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("E17").Value = "  -2.85%  "
$ws.Range("D18").Value = "'2.357.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.53%  "
$ws.Range("D19").Value = "'10.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.29%  "
$ws.Range("D20").Value = "'312.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.44%  "
$ws.Range("E21").Value = "  -4.30%  "
$ws.Range("D22").Value = "'6.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'65.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.13%  "
$ws.Range("D26").Value = "'2.499.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.17%  "
$ws.Range("D27").Value = "'0.378"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.15%  "
$ws.Range("E28").Value = "  -4.98%  "
$ws.Range("D29").Value = "'7.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").Value = "'174.83"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.18%  "
$ws.Range("E31").Value = "  -3.16%  "
$ws.Range("D32").Value = "'0.0₃0714"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("E34").Value = "  -5.96%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'0.994"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").Value = "'17.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("E39").Value = "  -4.36%  "
$ws.Range("D40").Value = "'35.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.38%  "
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("E42").Value = "  -5.74%  "
$ws.Range("D43").Value = "'132.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "
$ws.Range("E44").Value = "  -2.68%  "
$ws.Range("D45").Value = "'4.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.79%  "
$ws.Range("D46").Value = "'255.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.98%  "
$ws.Range("D47").Value = "'0.570"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.52%  "
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("E49").Value = "  -3.85%  "

# Rows 50 and 51 swapped their coin identity (EnergySwap <-> VeChain)
# along with updated price/volume figures
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0208"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.21%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'16.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.02%  "
